$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 8409.23  # H70: 8710 -> 8409.23
$ws.Cells.Item(70, 9).Value = 12505.111  # I70: 13957 -> 12505.111
$ws.Cells.Item(70, 11).Value = 37515.333  # K70: 41871 -> 37515.333
$ws.Cells.Item(70, 13).Value = -37245.333  # M70: -41601 -> -37245.333
$ws.Cells.Item(73, 8).Value = 8409.23  # H73: 8710 -> 8409.23
$ws.Cells.Item(73, 9).Value = 12505.111  # I73: 13957 -> 12505.111
$ws.Cells.Item(73, 11).Value = 37515.333  # K73: 41871 -> 37515.333
$ws.Cells.Item(73, 13).Value = -36579.333  # M73: -40935 -> -36579.333
$ws.Cells.Item(74, 8).Value = 2966.25  # H74: 2751.8572 -> 2966.25
$ws.Cells.Item(74, 9).Value = 2966.25  # I74: 2751.8572 -> 2966.25
$ws.Cells.Item(74, 11).Value = 2966.25  # K74: 2751.8572 -> 2966.25
$ws.Cells.Item(74, 13).Value = -2030.25  # M74: -1815.8572 -> -2030.25
$ws.Cells.Item(77, 8).Value = 2966.25  # H77: 2751.8572 -> 2966.25
$ws.Cells.Item(77, 9).Value = 2966.25  # I77: 2751.8572 -> 2966.25
$ws.Cells.Item(77, 11).Value = 14831.25  # K77: 13759.286 -> 14831.25
$ws.Cells.Item(77, 13).Value = -10151.25  # M77: -9079.286 -> -10151.25
$ws.Cells.Item(135, 8).Value = 15917.514  # H135: 16337.5 -> 15917.514
$ws.Cells.Item(135, 9).Value = 1670.8148  # I135: 1704.3846 -> 1670.8148
$ws.Cells.Item(135, 11).Value = 15037.3332  # K135: 15339.4614 -> 15037.3332
$ws.Cells.Item(135, 13).Value = -12502.3332  # M135: -12804.4614 -> -12502.3332

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 15517.134  # H32: 15517.256 -> 15517.134
$ws.Cells.Item(32, 9).Value = 15517.134  # I32: 15517.256 -> 15517.134
$ws.Cells.Item(32, 11).Value = 15517.134  # K32: 15517.256 -> 15517.134
$ws.Cells.Item(32, 13).Value = -15230.134  # M32: -15230.256 -> -15230.134
$ws.Cells.Item(45, 8).Value = 3308.8333  # H45: 3278.721 -> 3308.8333
$ws.Cells.Item(45, 10).Value = 4153.1  # J45: 4051.238 -> 4153.1
$ws.Cells.Item(45, 12).Value = 4153.1  # L45: 4051.238 -> 4153.1
$ws.Cells.Item(45, 14).Value = -4907.1  # N45: -4805.237999999999 -> -4907.1
$ws.Cells.Item(74, 8).Value = 215059.53  # H74: 207650.2 -> 215059.53
$ws.Cells.Item(74, 10).Value = 1200.25  # J74: 998 -> 1200.25
$ws.Cells.Item(74, 12).Value = 1200.25  # L74: 998 -> 1200.25
$ws.Cells.Item(74, 14).Value = -2948.25  # N74: -2746 -> -2948.25
$ws.Cells.Item(77, 8).Value = 215059.53  # H77: 207650.2 -> 215059.53
$ws.Cells.Item(77, 10).Value = 1200.25  # J77: 998 -> 1200.25
$ws.Cells.Item(77, 12).Value = 6001.25  # L77: 4990 -> 6001.25
$ws.Cells.Item(77, 14).Value = -14737.25  # N77: -13726 -> -14737.25
$ws.Cells.Item(101, 8).Value = 26999.5  # H101: 47860 -> 26999.5
$ws.Cells.Item(101, 10).Value = 26999.5  # J101: 47860 -> 26999.5
$ws.Cells.Item(101, 12).Value = 26999.5  # L101: 47860 -> 26999.5
$ws.Cells.Item(101, 14).Value = -33489.5  # N101: -54350 -> -33489.5
$ws.Cells.Item(132, 8).Value = 1756.0656  # H132: 1773.6833 -> 1756.0656
$ws.Cells.Item(132, 9).Value = 1272.122  # I132: 1286.45 -> 1272.122
$ws.Cells.Item(132, 11).Value = 3816.366  # K132: 3859.35 -> 3816.366
$ws.Cells.Item(132, 13).Value = -1286.366  # M132: -1329.35 -> -1286.366

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1340.9  # H86: 1431.3529 -> 1340.9
$ws.Cells.Item(86, 9).Value = 1441.75  # I86: 1527.3572 -> 1441.75
$ws.Cells.Item(86, 10).Value = 937.5  # J86: 983.3333 -> 937.5
$ws.Cells.Item(86, 11).Value = 1441.75  # K86: 1527.3572 -> 1441.75
$ws.Cells.Item(86, 12).Value = 937.5  # L86: 983.3333 -> 937.5
$ws.Cells.Item(86, 13).Value = -318.75  # M86: -404.3571999999999 -> -318.75
$ws.Cells.Item(86, 14).Value = -3183.5  # N86: -3229.3333 -> -3183.5
$ws.Cells.Item(89, 8).Value = 1340.9  # H89: 1431.3529 -> 1340.9
$ws.Cells.Item(89, 9).Value = 1441.75  # I89: 1527.3572 -> 1441.75
$ws.Cells.Item(89, 10).Value = 937.5  # J89: 983.3333 -> 937.5
$ws.Cells.Item(89, 11).Value = 7208.75  # K89: 7636.786 -> 7208.75
$ws.Cells.Item(89, 12).Value = 4687.5  # L89: 4916.6665 -> 4687.5
$ws.Cells.Item(89, 13).Value = -1592.75  # M89: -2020.786 -> -1592.75
$ws.Cells.Item(89, 14).Value = -15919.5  # N89: -16148.6665 -> -15919.5
$ws.Cells.Item(134, 8).Value = 2373.718  # H134: 2672.1177 -> 2373.718
$ws.Cells.Item(134, 9).Value = 2460.5151  # I134: 2659.3 -> 2460.5151
$ws.Cells.Item(134, 10).Value = 1896.3334  # J134: 2768.25 -> 1896.3334
$ws.Cells.Item(134, 11).Value = 7381.5453  # K134: 7977.900000000001 -> 7381.5453
$ws.Cells.Item(134, 12).Value = 5689.0002  # L134: 8304.75 -> 5689.0002
$ws.Cells.Item(134, 13).Value = -4846.5453  # M134: -5442.900000000001 -> -4846.5453
$ws.Cells.Item(134, 14).Value = -10759.0002  # N134: -13374.75 -> -10759.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 44.31579  # H7: 46 -> 44.31579
$ws.Cells.Item(7, 9).Value = 36.375  # I7: 37.214287 -> 36.375
$ws.Cells.Item(7, 10).Value = 86.666664  # J7: 87 -> 86.666664
$ws.Cells.Item(7, 11).Value = 36.375  # K7: 37.214287 -> 36.375
$ws.Cells.Item(7, 12).Value = 86.666664  # L7: 87 -> 86.666664
$ws.Cells.Item(7, 13).Value = 76.625  # M7: 75.785713 -> 76.625
$ws.Cells.Item(7, 14).Value = -312.666664  # N7: -313 -> -312.666664
$ws.Cells.Item(58, 8).Value = 1728.6471  # H58: 1792.7812 -> 1728.6471
$ws.Cells.Item(58, 9).Value = 1518.1482  # I58: 1579.12 -> 1518.1482
$ws.Cells.Item(58, 10).Value = 2540.5715  # J58: 2555.8572 -> 2540.5715
$ws.Cells.Item(58, 11).Value = 1518.1482  # K58: 1579.12 -> 1518.1482
$ws.Cells.Item(58, 12).Value = 2540.5715  # L58: 2555.8572 -> 2540.5715
$ws.Cells.Item(58, 13).Value = -1315.1482  # M58: -1376.12 -> -1315.1482
$ws.Cells.Item(58, 14).Value = -2946.5715  # N58: -2961.8572 -> -2946.5715
$ws.Cells.Item(59, 8).Value = 30011.2  # H59: 29514 -> 30011.2
$ws.Cells.Item(59, 9).Value = 24199.8  # I59: 24332.834 -> 24199.8
$ws.Cells.Item(59, 10).Value = 35822.6  # J59: 45057.5 -> 35822.6
$ws.Cells.Item(59, 11).Value = 24199.8  # K59: 24332.834 -> 24199.8
$ws.Cells.Item(59, 12).Value = 35822.6  # L59: 45057.5 -> 35822.6
$ws.Cells.Item(59, 13).Value = -23054.8  # M59: -23187.834 -> -23054.8
$ws.Cells.Item(59, 14).Value = -38112.6  # N59: -47347.5 -> -38112.6
$ws.Cells.Item(60, 8).Value = 36490.547  # H60: 37044.223 -> 36490.547
$ws.Cells.Item(60, 9).Value = 36899  # I60: 36799 -> 36899
$ws.Cells.Item(60, 10).Value = 36449.7  # J60: 37114.285 -> 36449.7
$ws.Cells.Item(60, 11).Value = 36899  # K60: 36799 -> 36899
$ws.Cells.Item(60, 12).Value = 36449.7  # L60: 37114.285 -> 36449.7
$ws.Cells.Item(60, 13).Value = -36388  # M60: -36288 -> -36388
$ws.Cells.Item(60, 14).Value = -37471.7  # N60: -38136.285 -> -37471.7
$ws.Cells.Item(134, 8).Value = 2805.3462  # H134: 2684.577 -> 2805.3462
$ws.Cells.Item(134, 9).Value = 1929.3684  # I134: 1897.9 -> 1929.3684
$ws.Cells.Item(134, 10).Value = 5183  # J134: 5306.8335 -> 5183
$ws.Cells.Item(134, 11).Value = 5788.1052  # K134: 5693.700000000001 -> 5788.1052
$ws.Cells.Item(134, 12).Value = 15549  # L134: 15920.5005 -> 15549
$ws.Cells.Item(134, 13).Value = -3253.1052  # M134: -3158.700000000001 -> -3253.1052
$ws.Cells.Item(134, 14).Value = -20619  # N134: -20990.5005 -> -20619
$ws.Cells.Item(136, 8).Value = 1728.6471  # H136: 1792.7812 -> 1728.6471
$ws.Cells.Item(136, 9).Value = 1518.1482  # I136: 1579.12 -> 1518.1482
$ws.Cells.Item(136, 10).Value = 2540.5715  # J136: 2555.8572 -> 2540.5715
$ws.Cells.Item(136, 11).Value = 4554.444600000001  # K136: 4737.36 -> 4554.444600000001
$ws.Cells.Item(136, 12).Value = 7621.7145  # L136: 7667.571599999999 -> 7621.7145
$ws.Cells.Item(136, 13).Value = -2004.444600000001  # M136: -2187.36 -> -2004.444600000001
$ws.Cells.Item(136, 14).Value = -12721.7145  # N136: -12767.5716 -> -12721.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 500003740  # H11: 1000000000 -> 500003740
$ws.Cells.Item(11, 9).Value = 666668350  # I11: 1000000000 -> 666668350
$ws.Cells.Item(11, 10).Value = 10004  # J11: 0 -> 10004
$ws.Cells.Item(11, 11).Value = 666668350  # K11: 1000000000 -> 666668350
$ws.Cells.Item(11, 12).Value = 10004  # L11: 0 -> 10004
$ws.Cells.Item(11, 13).Value = -666668211  # M11: -999999861 -> -666668211
$ws.Cells.Item(11, 14).Value = -10282  # N11: None -> -10282
$ws.Cells.Item(12, 8).Value = 15000  # H12: 0 -> 15000
$ws.Cells.Item(12, 10).Value = 15000  # J12: 0 -> 15000
$ws.Cells.Item(12, 12).Value = 15000  # L12: 0 -> 15000
$ws.Cells.Item(12, 14).Value = -15280  # N12: None -> -15280
$ws.Cells.Item(45, 8).Value = 49798.3  # H45: 49798.5 -> 49798.3
$ws.Cells.Item(45, 10).Value = 49798.3  # J45: 49798.5 -> 49798.3
$ws.Cells.Item(45, 12).Value = 49798.3  # L45: 49798.5 -> 49798.3
$ws.Cells.Item(45, 14).Value = -50916.3  # N45: -50916.5 -> -50916.3
$ws.Cells.Item(63, 8).Value = 29122.545  # H63: 30032.7 -> 29122.545
$ws.Cells.Item(63, 10).Value = 30027.223  # J63: 31278 -> 30027.223
$ws.Cells.Item(63, 12).Value = 30027.223  # L63: 31278 -> 30027.223
$ws.Cells.Item(63, 14).Value = -31399.223  # N63: -32650 -> -31399.223
$ws.Cells.Item(66, 8).Value = 29122.545  # H66: 30032.7 -> 29122.545
$ws.Cells.Item(66, 10).Value = 30027.223  # J66: 31278 -> 30027.223
$ws.Cells.Item(66, 12).Value = 90081.66900000001  # L66: 93834 -> 90081.66900000001
$ws.Cells.Item(66, 14).Value = -96945.66900000001  # N66: -100698 -> -96945.66900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 22004  # H17: 0 -> 22004
$ws.Cells.Item(17, 9).Value = 4008  # I17: 0 -> 4008
$ws.Cells.Item(17, 10).Value = 40000  # J17: 0 -> 40000
$ws.Cells.Item(17, 11).Value = 4008  # K17: 0 -> 4008
$ws.Cells.Item(17, 12).Value = 40000  # L17: 0 -> 40000
$ws.Cells.Item(17, 13).Value = -3838  # M17: None -> -3838
$ws.Cells.Item(17, 14).Value = -40340  # N17: None -> -40340
$ws.Cells.Item(22, 8).Value = 1431.2307  # H22: 1278.8096 -> 1431.2307
$ws.Cells.Item(22, 9).Value = 1001  # I22: 955 -> 1001
$ws.Cells.Item(22, 10).Value = 1622.4445  # J22: 1440.7142 -> 1622.4445
$ws.Cells.Item(22, 11).Value = 1001  # K22: 955 -> 1001
$ws.Cells.Item(22, 12).Value = 1622.4445  # L22: 1440.7142 -> 1622.4445
$ws.Cells.Item(22, 13).Value = -706  # M22: -660 -> -706
$ws.Cells.Item(22, 14).Value = -2212.4445  # N22: -2030.7142 -> -2212.4445
$ws.Cells.Item(27, 8).Value = 1431.2307  # H27: 1278.8096 -> 1431.2307
$ws.Cells.Item(27, 9).Value = 1001  # I27: 955 -> 1001
$ws.Cells.Item(27, 10).Value = 1622.4445  # J27: 1440.7142 -> 1622.4445
$ws.Cells.Item(27, 11).Value = 1001  # K27: 955 -> 1001
$ws.Cells.Item(27, 12).Value = 1622.4445  # L27: 1440.7142 -> 1622.4445
$ws.Cells.Item(27, 13).Value = -894  # M27: -848 -> -894
$ws.Cells.Item(27, 14).Value = -1836.4445  # N27: -1654.7142 -> -1836.4445
$ws.Cells.Item(40, 8).Value = 2152.75  # H40: 2226.8696 -> 2152.75
$ws.Cells.Item(40, 9).Value = 2141.3809  # I40: 2226.05 -> 2141.3809
$ws.Cells.Item(40, 11).Value = 2141.3809  # K40: 2226.05 -> 2141.3809
$ws.Cells.Item(40, 13).Value = -2005.3809  # M40: -2090.05 -> -2005.3809
$ws.Cells.Item(46, 8).Value = 3820.4  # H46: 4109.2607 -> 3820.4
$ws.Cells.Item(46, 10).Value = 5839  # J46: 6660.615 -> 5839
$ws.Cells.Item(46, 12).Value = 5839  # L46: 6660.615 -> 5839
$ws.Cells.Item(46, 14).Value = -6215  # N46: -7036.615 -> -6215
$ws.Cells.Item(93, 8).Value = 1670334.4  # H93: 1678168.1 -> 1670334.4
$ws.Cells.Item(93, 9).Value = 2503051.8  # I93: 2514802.5 -> 2503051.8
$ws.Cells.Item(93, 11).Value = 2503051.8  # K93: 2514802.5 -> 2503051.8
$ws.Cells.Item(93, 13).Value = -2501803.8  # M93: -2513554.5 -> -2501803.8
$ws.Cells.Item(94, 8).Value = 0  # H94: 30000 -> 0
$ws.Cells.Item(94, 10).Value = 0  # J94: 30000 -> 0
$ws.Cells.Item(94, 12).Value = 0  # L94: 30000 -> 0
$ws.Cells.Item(94, 14).ClearContents()  # N94: -31352 -> (removed)
$ws.Cells.Item(122, 8).Value = 7359.9067  # H122: 7359.93 -> 7359.9067
$ws.Cells.Item(122, 9).Value = 7839.143  # I122: 7839.1787 -> 7839.143
$ws.Cells.Item(122, 11).Value = 23517.429  # K122: 23517.5361 -> 23517.429
$ws.Cells.Item(122, 13).Value = -21067.429  # M122: -21067.5361 -> -21067.429
$ws.Cells.Item(132, 8).Value = 3539.3235  # H132: 3571.182 -> 3539.3235
$ws.Cells.Item(132, 9).Value = 3022.64  # I132: 3044.9167 -> 3022.64
$ws.Cells.Item(132, 11).Value = 9067.92  # K132: 9134.750100000001 -> 9067.92
$ws.Cells.Item(132, 13).Value = -6537.92  # M132: -6604.750100000001 -> -6537.92

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(33, 8).Value = 27755.25  # H33: 35510.5 -> 27755.25
$ws.Cells.Item(33, 10).Value = 27755.25  # J33: 35510.5 -> 27755.25
$ws.Cells.Item(33, 12).Value = 27755.25  # L33: 35510.5 -> 27755.25
$ws.Cells.Item(33, 14).Value = -28255.25  # N33: -36010.5 -> -28255.25
$ws.Cells.Item(36, 8).Value = 27755.25  # H36: 35510.5 -> 27755.25
$ws.Cells.Item(36, 10).Value = 27755.25  # J36: 35510.5 -> 27755.25
$ws.Cells.Item(36, 12).Value = 27755.25  # L36: 35510.5 -> 27755.25
$ws.Cells.Item(36, 14).Value = -28255.25  # N36: -36010.5 -> -28255.25
$ws.Cells.Item(64, 8).Value = 92996.664  # H64: 93500 -> 92996.664
$ws.Cells.Item(64, 9).Value = 91990  # I64: 0 -> 91990
$ws.Cells.Item(64, 11).Value = 91990  # K64: 0 -> 91990
$ws.Cells.Item(64, 13).Value = -91742  # M64: None -> -91742
$ws.Cells.Item(67, 8).Value = 92996.664  # H67: 93500 -> 92996.664
$ws.Cells.Item(67, 9).Value = 91990  # I67: 0 -> 91990
$ws.Cells.Item(67, 11).Value = 91990  # K67: 0 -> 91990
$ws.Cells.Item(67, 13).Value = -91132  # M67: None -> -91132
$ws.Cells.Item(107, 8).Value = 1207.8235  # H107: 1317.6 -> 1207.8235
$ws.Cells.Item(107, 9).Value = 1186.9166  # I107: 1347.4 -> 1186.9166
$ws.Cells.Item(107, 11).Value = 3560.7498  # K107: 4042.2 -> 3560.7498
$ws.Cells.Item(107, 13).Value = -1640.7498  # M107: -2122.2 -> -1640.7498
$ws.Cells.Item(122, 8).Value = 90761.80499999999  # H122: 98228.625 -> 90761.80499999999
$ws.Cells.Item(122, 9).Value = 106027.59  # I122: 116514.35 -> 106027.59
$ws.Cells.Item(122, 11).Value = 318082.77  # K122: 349543.05 -> 318082.77
$ws.Cells.Item(122, 13).Value = -315632.77  # M122: -347093.05 -> -315632.77
$ws.Cells.Item(126, 8).Value = 252327.5  # H126: 280291.66 -> 252327.5
$ws.Cells.Item(126, 9).Value = 2096.9375  # I126: 2303.6428 -> 2096.9375
$ws.Cells.Item(126, 11).Value = 6290.8125  # K126: 6910.928400000001 -> 6290.8125
$ws.Cells.Item(126, 13).Value = -3820.8125  # M126: -4440.928400000001 -> -3820.8125
$ws.Cells.Item(136, 8).Value = 17030.49  # H136: 17329.203 -> 17030.49
$ws.Cells.Item(136, 10).Value = 4848.9287  # J136: 5152.6924 -> 4848.9287
$ws.Cells.Item(136, 12).Value = 14546.7861  # L136: 15458.0772 -> 14546.7861
$ws.Cells.Item(136, 14).Value = -19646.7861  # N136: -20558.0772 -> -19646.7861
$ws.Cells.Item(138, 8).Value = 100000  # H138: 110000 -> 100000
$ws.Cells.Item(138, 10).Value = 100000  # J138: 110000 -> 100000
$ws.Cells.Item(138, 12).Value = 100000  # L138: 110000 -> 100000
$ws.Cells.Item(138, 14).Value = -110280  # N138: -120280 -> -110280
